$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 ("Rule" column for rule R40) is rewritten from the text "R40" to
# the text "1" (e.g. the rule is renumbered). The leading apostrophe forces
# Excel to keep storing a number-looking value as text, matching the
# original shared-string ("s") cell type instead of silently turning it
# into a numeric literal.
$ws.Range("B11").Value = "'1"
